# chore: update Sheets via scheduled runner
#
# Refreshes hard-coded marketboard price / profit figures on three of the
# "Rafflesia_Profits" Leve worksheets (ALC, BSM, WVR). These columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) are plain cached numbers, not formulas, so each cell is
# simply overwritten with its new value. Where the refreshed price feed no
# longer has data for a row, the previously-populated cells are cleared
# entirely (ClearContents) rather than zeroed, matching how the upstream
# job drops cells it has no data for.

$wb = $excel.ActiveWorkbook

# ---- ALC ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 33 ("Glazed and Confused"): price refresh, all columns stay populated.
$ws.Range("H33").Value = 753.3333
$ws.Range("I33").Value = 398.33334
$ws.Range("J33").Value = 1463.3334
$ws.Range("K33").Value = 398.33334
$ws.Range("L33").Value = 1463.3334
$ws.Range("M33").Value = -169.33334
$ws.Range("N33").Value = -1921.3334

# Row 103 ("Let Loose the Juice"): price feed now empty -> H:L go to 0,
# and the two profit columns (M, N) are dropped entirely.
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103:N103").ClearContents()

# ---- BSM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 86 ("Through Thick and Thin"): price refresh; HQ profit column (N)
# newly populated.
$ws.Range("H86").Value = 1307.6666
$ws.Range("I86").Value = 1409.5
$ws.Range("J86").Value = 798.5
$ws.Range("K86").Value = 1409.5
$ws.Range("L86").Value = 798.5
$ws.Range("M86").Value = -286.5
$ws.Range("N86").Value = -3044.5

# Row 89 ("Piercing Eyes Deserve Piercing Shafts (L)"): same pattern.
$ws.Range("H89").Value = 1307.6666
$ws.Range("I89").Value = 1409.5
$ws.Range("J89").Value = 798.5
$ws.Range("K89").Value = 7047.5
$ws.Range("L89").Value = 3992.5
$ws.Range("M89").Value = -1431.5
$ws.Range("N89").Value = -15224.5

# Rows 117-141: the refreshed run no longer carries price/profit data for
# these leves, so every previously-populated cell in H:N is cleared
# (column A-G identifiers/levels are untouched). Rows 121 and 136 already
# had no data beyond H:L=0 and are left as-is (not touched by the source
# diff).
$ws.Range("H117:M117").ClearContents()
$ws.Range("H118:M118").ClearContents()
$ws.Range("H119:L119").ClearContents()
$ws.Range("H120:L120").ClearContents()
$ws.Range("H122:L122").ClearContents()

$ws.Range("H123:L123").ClearContents()
$ws.Range("N123").ClearContents()

$ws.Range("H124:L124").ClearContents()
$ws.Range("N124").ClearContents()

$ws.Range("H125:L125").ClearContents()
$ws.Range("H126:L126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("H130:L130").ClearContents()

$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()

$ws.Range("H132:L132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H133:L133").ClearContents()
$ws.Range("H134:L134").ClearContents()

$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()

$ws.Range("H137:L137").ClearContents()

$ws.Range("H138:L138").ClearContents()
$ws.Range("N138").ClearContents()

$ws.Range("H139:L139").ClearContents()

$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()

$ws.Range("H141:L141").ClearContents()

# ---- WVR ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 100 ("Of Great Import"): price refresh (J/L/N are fixed leve
# rewards and stay unchanged).
$ws.Range("H100").Value = 1997.75
$ws.Range("I100").Value = 1830.3334
$ws.Range("K100").Value = 3660.6668
$ws.Range("M100").Value = -3119.6668

# Row 122 ("Heavy Armoire"): same pattern.
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 747.5
$ws.Range("K122").Value = 2242.5
$ws.Range("M122").Value = 207.5
